$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.774.56"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "2.307.66"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").Value = "2.308.66"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0955"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "2.691.60"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "54.591.24"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.286.23"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "308.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.376"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("E28").Value = "  +6.11%  "
$ws.Range("D29").Value = "2.399.88"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0700"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.377"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "128.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0896"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.552"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "245.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0489"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
